# Scrum Board update:
#  - "Começar a fazer o use case diagram" task is removed from the To-do
#    column (D4). The remaining To-do item ("Analisar o código dado",
#    previously in D5) moves up into D4, and D5 becomes empty again.
#  - The now-empty trailing cell D7 gets an underlined "blank slot" style.
#  - Selection cursor moves to D7.
#  - Page setup switched to A4 / portrait for printing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pull the next To-do item up into D4, then clear the now-vacated D5 cell.
$ws.Range("D4").Value = $ws.Range("D5").Value()
$ws.Range("D5").ClearContents()

# Give the new trailing blank cell an underlined style (matches the sheet's
# convention of underlining the next empty slot to fill in).
$ws.Range("D7").Font.Underline = $true

# Move the active selection to the newly underlined cell.
$ws.Range("D7").Select() | Out-Null

# Set print setup to A4, portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
